# The "Files" sheet lists one row per scanned DICOM file. A new "SessionLabel"
# value (column I) is being recorded for each file row, combining the
# subject ID with the scan/series identifier (e.g. "10001_CT1") - this is
# the work described in the commit message (collate_uploads computing a
# scan id instead of relying on DICOM:SeriesNumber alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")
$ws.Activate()

# Populate the new SessionLabel column (I) for the two data rows.
$ws.Range("I2").Value = "10001_CT1"
$ws.Range("I3").Value = "10001_CT1"

# Move the sheet's selection to the newly filled-in column, matching the
# cells the author was last working on.
$ws.Range("I2:I3").Select()
